$d = $word.ActiveDocument

# Merge split runs (that were separated due to grammar-check proofErr markers)
# back into a single run by performing a Find & Replace with the same final text.
# This causes Word to rewrite the matched range as one run, dropping the
# intervening w:proofErr elements.

$d.Content.Find.Execute(
    "Público-alvo: pessoas/alunos de 12 a 24 anos. Em caso do jogo ser usado de forma acadêmica, o público-alvo também inclui professores.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Público-alvo: pessoas/alunos de 12 a 24 anos. Em caso do jogo ser usado de forma acadêmica, o público-alvo também inclui professores.",
    2
)

$d.Content.Find.Execute(
    "Divertido: O jogo contará com um ranking dos melhores jogadores, com possibilidade deles ganharem cargos pelo seu desempenho no jogo.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Divertido: O jogo contará com um ranking dos melhores jogadores, com possibilidade deles ganharem cargos pelo seu desempenho no jogo.",
    2
)
